# Updates the cryptos list worksheet with refreshed price/volume data,
# and reorders the Celestia / TrustWalletToken rows (46 <-> 47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "44.188.99"
$ws.Range("E2").Value = "  +0.97%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "2.269.77"
$ws.Range("E3").Value = "  -0.12%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.04%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "231.08"

# --- Row 6 (XRP) ---
$ws.Range("E6").Value = "  +1.73%  "

# --- Row 7 (Solana) ---
$ws.Range("D7").Value = "64.38"
$ws.Range("E7").Value = "  +4.58%  "

# --- Row 8 (USDC) ---
$ws.Range("E8").Value = "  -0.05%  "

# --- Row 9 (Cardano) ---
# Leading apostrophe forces text so the trailing zero in "0.450" survives
# (a plain numeric-looking assignment would be coerced to 0.45).
$ws.Range("D9").Value = "'0.450"
$ws.Range("E9").Value = "  +5.88%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("D10").Value = "0.0999"
$ws.Range("E10").Value = "  +5.70%  "

# --- Row 11 (OKB) ---
$ws.Range("D11").Value = "57.22"
$ws.Range("E11").Value = "  -0.99%  "

# --- Row 12 (Avalanche) ---
$ws.Range("D12").Value = "27.54"
$ws.Range("E12").Value = "  +16.07%  "

# --- Row 13 (TRON) ---
$ws.Range("E13").Value = "  +1.80%  "

# --- Row 14 (WrappedliquidstakedEther2.0) ---
$ws.Range("D14").Value = "2.607.34"
$ws.Range("E14").Value = "  -0.13%  "

# --- Row 15 (Chainlink) ---
$ws.Range("D15").Value = "15.83"
$ws.Range("E15").Value = "  +0.32%  "

# --- Row 16 (Polkadot) ---
$ws.Range("D16").Value = "'6.10"
$ws.Range("E16").Value = "  +4.77%  "

# --- Row 17 (Polygon) ---
$ws.Range("D17").Value = "0.841"
$ws.Range("E17").Value = "  +3.67%  "

# --- Row 18 (WrappedEther) ---
$ws.Range("D18").Value = "2.282.99"
$ws.Range("E18").Value = "  +0.17%  "

# --- Row 19 (WrappedBTC) ---
$ws.Range("D19").Value = "44.020.51"
$ws.Range("E19").Value = "  +0.80%  "

# --- Row 20 (ShibaInu) ---
$ws.Range("E20").Value = "  +7.55%  "

# --- Row 21 (Litecoin) ---
$ws.Range("D21").Value = "73.89"
$ws.Range("E21").Value = "  +1.10%  "

# --- Row 22 (Uniswap) ---
$ws.Range("D22").Value = "6.12"
$ws.Range("E22").Value = "  -1.91%  "

# --- Row 23 (BitcoinCash) ---
$ws.Range("D23").Value = "253.44"
$ws.Range("E23").Value = "  +0.65%  "

# --- Row 24 (Dai) ---
$ws.Range("E24").Value = "  -0.03%  "

# --- Row 25 (PancakeSwap) ---
$ws.Range("E25").Value = "  -4.02%  "

# --- Row 26 (Cosmos) ---
$ws.Range("D26").Value = "10.15"
$ws.Range("E26").Value = "  +2.83%  "

# --- Row 27 (Toncoin) ---
$ws.Range("D27").Value = "2.25"
$ws.Range("E27").Value = "  -4.79%  "

# --- Row 28 (WEMIXToken) ---
$ws.Range("D28").Value = "3.24"
$ws.Range("E28").Value = "  +21.33%  "

# --- Row 29 (Monero) ---
$ws.Range("D29").Value = "171.58"
$ws.Range("E29").Value = "  +0.21%  "

# --- Row 30 (Kaspa) ---
$ws.Range("E30").Value = "  -0.02%  "

# --- Row 31 (EthereumClassic) ---
$ws.Range("D31").Value = "20.95"
$ws.Range("E31").Value = "  +1.67%  "

# --- Row 32 (ImmutableX) ---
$ws.Range("E32").Value = "  -4.12%  "

# --- Row 33 (Stellar) ---
$ws.Range("E33").Value = "  +2.78%  "

# --- Row 34 (Hedera) ---
$ws.Range("D34").Value = "0.0707"
$ws.Range("E34").Value = "  +6.49%  "

# --- Row 35 (Filecoin) ---
$ws.Range("D35").Value = "4.81"
$ws.Range("E35").Value = "  +0.62%  "

# --- Row 36 (InternetComputer(DFINITY)) ---
$ws.Range("D36").Value = "'4.90"
$ws.Range("E36").Value = "  -3.48%  "

# --- Row 37 (RenderToken) ---
$ws.Range("D37").Value = "3.82"
$ws.Range("E37").Value = "  +5.59%  "

# --- Row 38 (THORChain) ---
$ws.Range("E38").Value = "  +0.48%  "

# --- Row 39 (LidoDAOToken) ---
$ws.Range("E39").Value = "  -3.80%  "

# --- Row 40 (VeChain) ---
$ws.Range("D40").Value = "'0.0260"
$ws.Range("E40").Value = "  +3.81%  "

# --- Row 41 (BinanceUSD) ---
$ws.Range("E41").Value = "  -0.09%  "

# --- Row 42 (TerraClassic) ---
$ws.Range("E42").Value = "  -0.72%  "

# --- Row 43 (Cronos) ---
$ws.Range("D43").Value = "0.0991"
$ws.Range("E43").Value = "  +0.36%  "

# --- Row 44 (InjectiveProtocol) ---
$ws.Range("D44").Value = "17.52"
$ws.Range("E44").Value = "  +5.26%  "

# --- Row 45 (FraxShare) ---
$ws.Range("E45").Value = "  -6.13%  "

# --- Rows 46/47: Celestia and TrustWalletToken swap places ---
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "1.21"
$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").Value = "10.43"
$ws.Range("E47").Value = "  +6.86%  "

# --- Row 48 (Aave) ---
$ws.Range("D48").Value = "'98.50"
$ws.Range("E48").Value = "  +0.42%  "

# --- Row 49 (FTXToken) ---
$ws.Range("E49").Value = "  -2.56%  "

# --- Row 50 (NEARProtocol) ---
$ws.Range("D50").Value = "2.35"
$ws.Range("E50").Value = "  +2.99%  "

# --- Row 51 (Maker) ---
$ws.Range("D51").Value = "1.447.69"
$ws.Range("E51").Value = "  -1.86%  "
